$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "G8"  = 0
    "G11" = 3.38
    "G12" = 0.15
    "G13" = 0.14
    "G14" = 0.14
    "G16" = 0.89
    "G17" = 0.93
    "G18" = 0.16
    "G19" = 0.51
    "G20" = 6.53
    "G21" = 3.68
    "G22" = 3.62
    "G23" = 3.56
    "G24" = 1.03
    "G25" = 6.32
    "G26" = 6.51
    "G27" = 3.48
    "G28" = 5.79
    "G29" = 0.14
    "G32" = 0.15
    "G33" = 0.17
    "G38" = 0.29
    "G39" = 0.3
    "G41" = 0.32
    "G42" = 0.31
    "G46" = 0.29
    "G49" = 0.44
    "G50" = 0.44
    "G51" = 0.42
    "G52" = 0.51
    "G54" = 0.46
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
